$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.278.16"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "1.863.53"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'235.65"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").Value = "'0.2838"
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("D9").Value = "'0.06532"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "'21.72"
$ws.Range("E10").Value = "  +9.00%  "
$ws.Range("D11").Value = "'0.07940"
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").Value = "1.870.03"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").Value = "'5.162"
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("D15").Value = "'0.6787"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "'279.60"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("D17").Value = "30.284.62"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "'13.40"
$ws.Range("E18").Value = "  +6.68%  "
$ws.Range("D20").Value = "'5.395"
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("D21").Value = "2.113.31"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").Value = "'0.000007310"
$ws.Range("E22").Value = "  +1.27%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'6.171"
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("D25").Value = "'167.48"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").Value = "'9.176"
$ws.Range("E26").Value = "  -1.07%  "
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("E29").Value = "  +3.27%  "
$ws.Range("D30").Value = "'0.09745"
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("D31").Value = "'4.395"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("D33").Value = "'4.071"
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("D34").Value = "'0.04738"
$ws.Range("E34").Value = "  +1.50%  "
$ws.Range("E35").Value = "  +3.57%  "
$ws.Range("D36").Value = "'0.7063"
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("D37").Value = "'2.709"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("D39").Value = "'2.580"
$ws.Range("E39").Value = "  +2.88%  "
$ws.Range("D40").Value = "'6.319"
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("D41").Value = "'75.03"
$ws.Range("E41").Value = "  +4.66%  "
$ws.Range("D42").Value = "'1.961"
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("D43").Value = "'0.8508"
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("D44").Value = "'0.4183"
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'103.34"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").Value = "'965.61"
$ws.Range("E47").Value = "  -5.67%  "
$ws.Range("D48").Value = "'7.188"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").Value = "'9.343"
$ws.Range("E49").Value = "  +4.80%  "
$ws.Range("D50").Value = "'34.10"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("D51").Value = "'0.1134"
$ws.Range("E51").Value = "  -0.43%  "
